# Daily attendance processing - normalize the "Recorded By" (column G) list
# so the comma-separated names/emails are in case-sensitive (ordinal) sorted
# order, matching the canonical formatting used by the attendance system.

function Compare-Ordinal($s1, $s2) {
    $len1 = $s1.Length
    $len2 = $s2.Length
    $minLen = [Math]::Min($len1, $len2)
    for ($k = 0; $k -lt $minLen; $k++) {
        $c1 = [int][char]$s1[$k]
        $c2 = [int][char]$s2[$k]
        if ($c1 -lt $c2) { return -1 }
        if ($c1 -gt $c2) { return 1 }
    }
    if ($len1 -lt $len2) { return -1 }
    if ($len1 -gt $len2) { return 1 }
    return 0
}

function Sort-Ordinal($items) {
    $n = $items.Count
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt ($n - $i - 1); $j++) {
            $cmp = Compare-Ordinal $items[$j] $items[$j + 1]
            if ($cmp -gt 0) {
                $tmp = $items[$j]
                $items[$j] = $items[$j + 1]
                $items[$j + 1] = $tmp
            }
        }
    }
    return $items
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = @($val -split ", ")
        if ($parts.Count -gt 1) {
            $sorted = Sort-Ordinal $parts
            $newVal = $sorted -join ", "
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
